# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# The workbook is a "metadata" sheet describing dimensions/measures for a
# statistical table. Re-curation swaps which columns are treated as
# curated dimensions vs measures:
#   - "municipio-nombre" (column C) moves from measure -> dimension
#     (now uses sdmx-dimension:refArea / dim / URI-Municipio, matching the
#     other refArea columns D, J).
#   - "porcentaje-sau-regimen-tenencia" (column E) moves from dimension ->
#     measure (iaest-measure:... / medida / xsd:int), so its mapping-file
#     reference in row 5 is no longer needed and is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("municipio-nombre"): measure -> dimension
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("C3").Value = "dim"
$ws.Range("C4").Value = "URI-Municipio"

# Column E ("porcentaje-sau-regimen-tenencia"): dimension -> measure
$ws.Range("E2").Value = "iaest-measure:porcentaje-sau-regimen-tenencia"
$ws.Range("E3").Value = "medida"
$ws.Range("E4").Value = "xsd:int"

# Row 5: the porcentaje-sau-regimen-tenencia column no longer needs its
# curated mapping-workbook reference. Clear() (not just ClearContents)
# drops the cell record entirely, same as the other untouched cells
# (B5/C5/D5) in that row which were never populated.
$ws.Range("E5").Clear()
